# "contingencies with rene fine"
# Build a small 2x2 block:
#   B1 = 0              (bold, centered, thin box border)
#   A2 = 0              (bold, centered, thin box border - same style as B1)
#   B2 = "disconnected_elements" (plain text, default formatting)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values first.
$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

# Build the full style (bold font + thin box border + centered/top aligned)
# on B1 first.
$r1 = $ws.Range("B1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment = -4160     # xlTop
$r1.Borders.LineStyle = 1         # xlContinuous
$r1.Borders.Weight = 2            # xlThin

# Copy B1's formatting (only) onto A2 so both cells collapse onto the exact
# same cell-style record instead of each accumulating its own (which is
# what happens if the same property sequence is replayed on a second,
# independent range).
$r1.Copy()
$r2 = $ws.Range("A2")
$r2.PasteSpecial(-4122)           # xlPasteFormats
$excel.CutCopyMode = $false
